# Sprint 2 stories update
# - Mark two Backlog items (US35 "List recent births", US36 "List recent deaths")
#   as assigned to Sprint 2 / Planned.
# - Populate the Sprint2 worksheet with those two stories.
# - Restore view/selection state to match what was left after the edits
#   (Sprint2 becomes the active/selected sheet).

$wb = $excel.ActiveWorkbook

# ---- Backlog: rows 7 & 8 move from "Not Planned" to Sprint 2 / "Planned" ----
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("A7").Value = 2
$backlog.Range("E7").Value = "Planned"
$backlog.Range("A8").Value = 2
$backlog.Range("E8").Value = "Planned"

# ---- Sprint2: fill in the two stories pulled from the backlog ----
$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Range("A2").Value = "US35"
$sprint2.Range("B2").Value = "List recent births"
$sprint2.Range("C2").Value = "DA"
$sprint2.Range("D2").Value = "Planned"
$sprint2.Range("E2").Value = 25
$sprint2.Range("F2").Value = 60

$sprint2.Range("A3").Value = "US36"
$sprint2.Range("B3").Value = "List recent deaths"
$sprint2.Range("C3").Value = "DA"
$sprint2.Range("D3").Value = "Planned"
$sprint2.Range("E3").Value = 25
$sprint2.Range("F3").Value = 60

# ---- View/selection bookkeeping to mirror the saved workbook state ----
$backlog.Activate()
$backlog.Range("F26").Select()

$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Activate()
$sprint1.Range("C20").Select()

$stories = $wb.Worksheets.Item("Stories")
$stories.Activate()
$stories.Range("C37").Select()

# Sprint2 is the last activated / active sheet when the file was saved.
$sprint2.Activate()
$sprint2.Range("F9").Select()
